$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = -12.406
$ws.Range("B3").Value = 6.247999999999999
$ws.Range("C5").Value = -13.066
$ws.Range("D5").Value = -8.021000000000001
$ws.Range("D9").Value = -7.902000000000001
$ws.Range("D11").Value = -8.021000000000001
$ws.Range("B14").Value = 6.449
$ws.Range("B16").Value = 6.016
$ws.Range("C16").Value = -11.856
$ws.Range("D17").Value = -7.784000000000001
$ws.Range("B21").Value = 6.054
$ws.Range("D21").Value = -7.571
$ws.Range("B23").Value = 6.444
$ws.Range("B25").Value = 6.103999999999999
